# Updated simulated-game transition-probability matrix (more games simulated,
# simulate-game logic sped up, optimization logic drafted).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2192982456140351
$ws.Range("C2").Value = 0.5219298245614035
$ws.Range("J2").Value = 0.0131578947368421
$ws.Range("P2").Value = 0.1403508771929824
$ws.Range("S2").Value = 0.1052631578947368

# Row 3
$ws.Range("C3").Value = 0.008264462809917356
$ws.Range("J3").Value = 0.01652892561983471
$ws.Range("P3").Value = 0.7851239669421488
$ws.Range("S3").Value = 0.1900826446280992

# Row 4
$ws.Range("J4").Value = 0.0196078431372549
$ws.Range("P4").Value = 0.7450980392156863
$ws.Range("S4").Value = 0.2352941176470588

# Row 6
$ws.Range("B6").Value = 0.0427807486631016
$ws.Range("D6").Value = 0.0053475935828877
$ws.Range("F6").Value = 0.0427807486631016
$ws.Range("J6").Value = 0.2887700534759358
$ws.Range("O6").Value = 0.0053475935828877
$ws.Range("Q6").Value = 0.1711229946524064
$ws.Range("R6").Value = 0.053475935828877
$ws.Range("S6").Value = 0.3903743315508021

# Row 7
$ws.Range("B7").Value = 0.1038961038961039
$ws.Range("D7").Value = 0.04545454545454546
$ws.Range("F7").Value = 0.06493506493506493
$ws.Range("J7").Value = 0.1233766233766234
$ws.Range("O7").Value = 0.01948051948051948
$ws.Range("Q7").Value = 0.1883116883116883
$ws.Range("R7").Value = 0.07792207792207792
$ws.Range("S7").Value = 0.3766233766233766

# Row 8
$ws.Range("B8").Value = 0.08611111111111111
$ws.Range("D8").Value = 0.01666666666666667
$ws.Range("F8").Value = 0.06111111111111111
$ws.Range("J8").Value = 0.1138888888888889
$ws.Range("O8").Value = 0.01388888888888889
$ws.Range("Q8").Value = 0.1777777777777778
$ws.Range("R8").Value = 0.1027777777777778
$ws.Range("S8").Value = 0.4277777777777778

# Row 9
$ws.Range("B9").Value = 0.08609271523178808
$ws.Range("D9").Value = 0.006622516556291391
$ws.Range("F9").Value = 0.0728476821192053
$ws.Range("J9").Value = 0.09271523178807947
$ws.Range("O9").Value = 0.01324503311258278
$ws.Range("Q9").Value = 0.2317880794701987
$ws.Range("R9").Value = 0.06622516556291391
$ws.Range("S9").Value = 0.4304635761589404

# Row 10
$ws.Range("B10").Value = 0.1084812623274162
$ws.Range("D10").Value = 0.03353057199211045
$ws.Range("E10").Value = 0.0009861932938856016
$ws.Range("F10").Value = 0.0650887573964497
$ws.Range("J10").Value = 0.1252465483234714
$ws.Range("O10").Value = 0.01676528599605523
$ws.Range("Q10").Value = 0.2495069033530572
$ws.Range("R10").Value = 0.07692307692307693
$ws.Range("S10").Value = 0.3234714003944773

# Row 11
$ws.Range("G11").Value = 0.1111111111111111
$ws.Range("J11").Value = 0.09259259259259259
$ws.Range("K11").Value = 0.162037037037037
$ws.Range("L11").Value = 0.6111111111111112
$ws.Range("S11").Value = 0.02314814814814815

# Row 12
$ws.Range("G12").Value = 0.8029197080291971
$ws.Range("J12").Value = 0.1386861313868613
$ws.Range("K12").Value = 0.0145985401459854
$ws.Range("L12").Value = 0.0218978102189781
$ws.Range("S12").Value = 0.0218978102189781

# Row 13
$ws.Range("G13").Value = 0.7297297297297297
$ws.Range("J13").Value = 0.1891891891891892
$ws.Range("S13").Value = 0.08108108108108109

# Row 15
$ws.Range("F15").Value = 0.01578947368421053
$ws.Range("H15").Value = 0.1526315789473684
$ws.Range("I15").Value = 0.05789473684210526
$ws.Range("J15").Value = 0.4
$ws.Range("K15").Value = 0.03684210526315789
$ws.Range("M15").Value = 0.01052631578947368
$ws.Range("O15").Value = 0.05263157894736842
$ws.Range("S15").Value = 0.2736842105263158

# Row 16
$ws.Range("F16").Value = 0.03726708074534162
$ws.Range("H16").Value = 0.1366459627329193
$ws.Range("I16").Value = 0.1055900621118012
$ws.Range("J16").Value = 0.3664596273291926
$ws.Range("K16").Value = 0.09316770186335403
$ws.Range("M16").Value = 0.01863354037267081
$ws.Range("N16").Value = 0.006211180124223602
$ws.Range("O16").Value = 0.05590062111801242
$ws.Range("S16").Value = 0.1801242236024845

# Row 17
$ws.Range("F17").Value = 0.01951219512195122
$ws.Range("H17").Value = 0.1853658536585366
$ws.Range("I17").Value = 0.08536585365853659
$ws.Range("J17").Value = 0.4097560975609756
$ws.Range("K17").Value = 0.07804878048780488
$ws.Range("M17").Value = 0.02195121951219512
$ws.Range("O17").Value = 0.07804878048780488
$ws.Range("S17").Value = 0.1219512195121951

# Row 18
$ws.Range("F18").Value = 0.04137931034482759
$ws.Range("H18").Value = 0.2068965517241379
$ws.Range("I18").Value = 0.1103448275862069
$ws.Range("J18").Value = 0.3793103448275862
$ws.Range("K18").Value = 0.0896551724137931
$ws.Range("M18").Value = 0.01379310344827586
$ws.Range("O18").Value = 0.04827586206896552
$ws.Range("S18").Value = 0.1103448275862069

# Row 19
$ws.Range("F19").Value = 0.01814516129032258
$ws.Range("H19").Value = 0.2066532258064516
$ws.Range("I19").Value = 0.07056451612903226
$ws.Range("J19").Value = 0.3618951612903226
$ws.Range("K19").Value = 0.1108870967741935
$ws.Range("M19").Value = 0.02318548387096774
$ws.Range("O19").Value = 0.1169354838709677
$ws.Range("S19").Value = 0.09173387096774194
